$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new report row (row 17) with the LoadBalancedReport entry
$ws.Range("A17").Value = "Report_KHSH_GetLPPostLoadBalance"
$ws.Range("B17").Value = "2021 May 30"

# "1400/03/09" looks like a date, so a direct .Value assignment would get
# auto-converted into a date serial number. Build it as literal text in a
# scratch cell via a formula, then copy/paste-special (values only) into
# C17 so it lands as plain text without altering the cell's existing style.
$ws.Range("Z1").Formula = "=""1400/03/09"""
$ws.Range("Z1").Copy()
$ws.Range("C17").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# Update the active selection to reflect where the author left off
$ws.Range("C18").Select()
